# Appends 53 new bank-transaction rows (2023-09-05 .. 2023-11-14) to the "pcbanking"
# worksheet, continuing the existing ledger that ends at row 2425 (date 2023-08-31).
# Each new row reuses the date-format style of the preceding row (style index 2,
# format "m/d/yyyy") by copying that cell's formatting before writing the new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 2426
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/5/2023"
$ws.Range("B" + $r).Value = -29.89
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Miscellaneous Payment  "
$ws.Range("E" + $r).Value = "PAYPAL                       "

$r = 2427
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/5/2023"
$ws.Range("B" + $r).Value = -500
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Bill Payment           "
$ws.Range("E" + $r).Value = "PC-SCOTIABANK GOLD AMEX CARD"

$r = 2428
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/7/2023"
$ws.Range("B" + $r).Value = -134.25
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "CRD. Card Bill Payment "

$r = 2429
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/7/2023"
$ws.Range("B" + $r).Value = 1470.33
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Payroll Deposit        "
$ws.Range("E" + $r).Value = "BWS MANUFACTURING LTD        "

$r = 2430
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/12/2023"
$ws.Range("B" + $r).Value = -230
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "CNB #021                 WOODS"

$r = 2431
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/12/2023"
$ws.Range("B" + $r).Value = -157
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS CNB #021            WOODS"

$r = 2432
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/20/2023"
$ws.Range("B" + $r).Value = -2.36
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS SH VENDING8003620026MIRAM"

$r = 2433
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/20/2023"
$ws.Range("B" + $r).Value = -2.95
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS SH VENDING8003620026MIRAM"

$r = 2434
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/20/2023"
$ws.Range("B" + $r).Value = -9.19
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS FLORENCEVILLE IRVINGWEST "

$r = 2435
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/21/2023"
$ws.Range("B" + $r).Value = -34.49
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "CENTREVILLE VALU FOODS   CENTR"

$r = 2436
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/21/2023"
$ws.Range("B" + $r).Value = 1470.33
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Payroll Deposit        "
$ws.Range("E" + $r).Value = "BWS MANUFACTURING LTD        "

$r = 2437
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/25/2023"
$ws.Range("B" + $r).Value = -400
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Bill Payment           "
$ws.Range("E" + $r).Value = "PC-SCOTIABANK GOLD AMEX CARD"

$r = 2438
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/29/2023"
$ws.Range("B" + $r).Value = -34.49
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS CENTREVILLE VALU FOOCENTR"

$r = 2439
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "9/29/2023"
$ws.Range("B" + $r).Value = -9.19
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS FLORENCEVILLE IRVINGWEST "

$r = 2440
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/3/2023"
$ws.Range("B" + $r).Value = -600
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Miscellaneous Payment  "
$ws.Range("E" + $r).Value = "NSLSC                        "

$r = 2441
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/3/2023"
$ws.Range("B" + $r).Value = -1000
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Bill Payment           "
$ws.Range("E" + $r).Value = "PC-SCOTIABANK GOLD AMEX CARD"

$r = 2442
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/3/2023"
$ws.Range("B" + $r).Value = -1000
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Bill Payment           "
$ws.Range("E" + $r).Value = "PC-SCOTIABANK GOLD AMEX CARD"

$r = 2443
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/4/2023"
$ws.Range("B" + $r).Value = -12.52
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS MCDONALD'S #17867   WOODS"

$r = 2444
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/4/2023"
$ws.Range("B" + $r).Value = -135
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS CNB #021            WOODS"

$r = 2445
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/4/2023"
$ws.Range("B" + $r).Value = -215
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "CNB #021                 WOODS"

$r = 2446
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/5/2023"
$ws.Range("B" + $r).Value = 106.25
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "GST                    "
$ws.Range("E" + $r).Value = "CANADA                       "

$r = 2447
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/5/2023"
$ws.Range("B" + $r).Value = 1470.33
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Payroll Deposit        "
$ws.Range("E" + $r).Value = "BWS MANUFACTURING LTD        "

$r = 2448
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/6/2023"
$ws.Range("B" + $r).Value = -34.49
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS CENTREVILLE VALU FOOCENTR"

$r = 2449
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/11/2023"
$ws.Range("B" + $r).Value = -5.27
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS ESSO COUCHE-TARD    WEST "

$r = 2450
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/11/2023"
$ws.Range("B" + $r).Value = -750
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Bill Payment           "
$ws.Range("E" + $r).Value = "PC-SCOTIABANK GOLD AMEX CARD"

$r = 2451
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/13/2023"
$ws.Range("B" + $r).Value = 184
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Climate Action Incentive"
$ws.Range("E" + $r).Value = "CANADA                       "

$r = 2452
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/14/2023"
$ws.Range("B" + $r).Value = -2.07
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS MCDONALD'S #29096   FREDE"

$r = 2453
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/14/2023"
$ws.Range("B" + $r).Value = -2.07
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS TIM HORTONS #0423   FREDE"

$r = 2454
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/19/2023"
$ws.Range("B" + $r).Value = 1470.33
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Payroll Deposit        "
$ws.Range("E" + $r).Value = "BWS MANUFACTURING LTD        "

$r = 2455
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/19/2023"
$ws.Range("B" + $r).Value = -1000
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Bill Payment           "
$ws.Range("E" + $r).Value = "PC-SCOTIABANK GOLD AMEX CARD"

$r = 2456
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/19/2023"
$ws.Range("B" + $r).Value = -8000
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "RIVERVIEW HONDA          FREDE"

$r = 2457
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/21/2023"
$ws.Range("B" + $r).Value = -78.84
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "WITHDRAWAL             "
$ws.Range("E" + $r).Value = "FREE INTERAC E-TRANSFER"

$r = 2458
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/26/2023"
$ws.Range("B" + $r).Value = -21.43
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "GCDS625 BRISTOL          BRIST"

$r = 2459
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/26/2023"
$ws.Range("B" + $r).Value = -100.18
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "WALMART STORE #1043      WOODS"

$r = 2460
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/26/2023"
$ws.Range("B" + $r).Value = -19.18
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "DOLLARAMA # 167          WOODS"

$r = 2461
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/26/2023"
$ws.Range("B" + $r).Value = -155
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "CNB #021                 WOODS"

$r = 2462
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/26/2023"
$ws.Range("B" + $r).Value = -200
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "CNB #021                 WOODS"

$r = 2463
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/27/2023"
$ws.Range("B" + $r).Value = -2.06
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS SH VENDING8003620026MIRAM"

$r = 2464
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/28/2023"
$ws.Range("B" + $r).Value = -13.25
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS CENTRE BELL CONCESSIMONTR"

$r = 2465
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/28/2023"
$ws.Range("B" + $r).Value = -29.15
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS CENTRE BELL CONCESSIMONTR"

$r = 2466
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/28/2023"
$ws.Range("B" + $r).Value = -7.75
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS STADE SAPUTO - CONCEMONTR"

$r = 2467
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/28/2023"
$ws.Range("B" + $r).Value = -30
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS PARC OLYMPIQUE STATIMONTR"

$r = 2468
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/30/2023"
$ws.Range("B" + $r).Value = -8.1
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS DEPANNEUR P.BEDARD &VILLE"

$r = 2469
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "10/31/2023"
$ws.Range("B" + $r).Value = -600
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Miscellaneous Payment  "
$ws.Range("E" + $r).Value = "NSLSC                        "

$r = 2470
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "11/2/2023"
$ws.Range("B" + $r).Value = -73.59
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Miscellaneous Payment  "
$ws.Range("E" + $r).Value = "PAYPAL                       "

$r = 2471
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "11/2/2023"
$ws.Range("B" + $r).Value = 1470.33
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Payroll Deposit        "
$ws.Range("E" + $r).Value = "BWS MANUFACTURING LTD        "

$r = 2472
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "11/3/2023"
$ws.Range("B" + $r).Value = -366.69
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Loans                  "
$ws.Range("E" + $r).Value = "HONDA CANADA FINANCE INC.    "

$r = 2473
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "11/3/2023"
$ws.Range("B" + $r).Value = -214
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "MARK'S THE SPOT          BEECH"

$r = 2474
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "11/3/2023"
$ws.Range("B" + $r).Value = -3000
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Bill Payment           "
$ws.Range("E" + $r).Value = "PC-SCOTIABANK GOLD AMEX CARD"

$r = 2475
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "11/8/2023"
$ws.Range("B" + $r).Value = -780.37
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Bill Payment           "
$ws.Range("E" + $r).Value = "PC-SCOTIABANK GOLD AMEX CARD"

$r = 2476
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "11/10/2023"
$ws.Range("B" + $r).Value = -10.4
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "Miscellaneous Payment  "
$ws.Range("E" + $r).Value = "PAYPAL                       "

$r = 2477
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "11/14/2023"
$ws.Range("B" + $r).Value = -8.03
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS TIM HORTONS #2853   FLORE"

$r = 2478
$ws.Range("A" + ($r - 1)).Copy()
$ws.Range("A" + $r).PasteSpecial(-4122)
$ws.Range("A" + $r).Value = "11/14/2023"
$ws.Range("B" + $r).Value = -20.37
$ws.Range("C" + $r).Value = "-"
$ws.Range("D" + $r).Value = "POS Purchase           "
$ws.Range("E" + $r).Value = "FPOS CO-OP #9283         FLORE"

$excel.CutCopyMode = $false

# Match the author's final cursor position/selection in the sheet view.
$ws.Range("D2422").Select() | Out-Null
